$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay stored as literal text
# (matches the source data which keeps things like trailing zeros, e.g. "0.5400")
$ws.Range("D2").Value = "26.136.68"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.667.18"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.89"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5208"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2602"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06328"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.06"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07533"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.661.27"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.406"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5400"
$ws.Range("E14").Value = "  -5.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000007984"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.36"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "26.174.13"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.724"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.97"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("E21").Value = "  -3.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.222"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.62"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1235"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.407"
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.72"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06267"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.274"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.489"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.403"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.633"
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9977"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.764"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.393"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5966"
$ws.Range("E37").Value = "  -1.93%  "
$ws.Range("D38").Value = "1.108.75"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.052"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8558"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.67"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "1.815.30"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.31"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.029"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05238"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4238"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.877"
$ws.Range("E51").Value = "  -1.65%  "
